# Applies the "Updated cryptos list" data refresh to Sheet1 (rows 2-51).
# Column A (rank index) is untouched; B/C only change where coins were
# re-ordered (rows 20/21 and 44/45 swap places); D/E are refreshed for every row.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.712.72"
$ws.Range("E2").Value = "  -1.68%  "

$ws.Range("D3").Value = "2.553.66"
$ws.Range("E3").Value = "  -4.80%  "

$ws.Range("D4").Value = "'1.00"
$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'577.93"
$ws.Range("E5").Value = "  -3.46%  "

$ws.Range("D6").Value = "'170.42"
$ws.Range("E6").Value = "  -2.72%  "

$ws.Range("D7").Value = "'1.00"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").Value = "'0.511"
$ws.Range("E8").Value = "  -2.28%  "

$ws.Range("D9").Value = "2.554.51"
$ws.Range("E9").Value = "  -4.68%  "

$ws.Range("D10").Value = "'0.166"
$ws.Range("E10").Value = "  -0.66%  "

$ws.Range("D11").Value = "'0.169"
$ws.Range("E11").Value = "  -0.11%  "

$ws.Range("D12").Value = "'0.346"
$ws.Range("E12").Value = "  -2.79%  "

$ws.Range("D13").Value = "'4.82"
$ws.Range("E13").Value = "  -3.16%  "

$ws.Range("D14").Value = "3.023.03"
$ws.Range("E14").Value = "  -4.65%  "

$ws.Range("D15").Value = "'0.0000181"
$ws.Range("E15").Value = "  -1.88%  "

$ws.Range("D16").Value = "70.542.71"
$ws.Range("E16").Value = "  -1.78%  "

$ws.Range("D17").Value = "'25.05"
$ws.Range("E17").Value = "  -4.39%  "

$ws.Range("D18").Value = "2.564.69"
$ws.Range("E18").Value = "  -4.29%  "

$ws.Range("D19").Value = "'11.66"
$ws.Range("E19").Value = "  -4.39%  "

$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").Value = "'362.87"
$ws.Range("E20").Value = "  -2.26%  "

$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D21").Value = "'7.09"
$ws.Range("E21").Value = "  -13.87%  "

$ws.Range("D22").Value = "'3.96"
$ws.Range("E22").Value = "  -5.19%  "

$ws.Range("D23").Value = "'2.00"
$ws.Range("E23").Value = "  -0.58%  "

$ws.Range("D24").Value = "'1.00"
$ws.Range("E24").Value = "  +0.09%  "

$ws.Range("D25").Value = "'70.05"
$ws.Range("E25").Value = "  -2.84%  "

$ws.Range("D26").Value = "'4.11"
$ws.Range("E26").Value = "  -5.09%  "

$ws.Range("D27").Value = "'9.31"
$ws.Range("E27").Value = "  -4.43%  "

$ws.Range("D28").Value = "2.687.13"
$ws.Range("E28").Value = "  -4.60%  "

$ws.Range("D29").Value = "'0.998"
$ws.Range("E29").Value = "  -0.06%  "

$ws.Range("D30").Value = "0.0⃳0930"
$ws.Range("E30").Value = "  -4.17%  "

$ws.Range("D31").Value = "'7.90"
$ws.Range("E31").Value = "  -1.80%  "

$ws.Range("D32").Value = "'487.48"
$ws.Range("E32").Value = "  -3.06%  "

$ws.Range("D33").Value = "'1.29"
$ws.Range("E33").Value = "  -0.58%  "

$ws.Range("D34").Value = "'1.77"
$ws.Range("E34").Value = "  -2.85%  "

$ws.Range("D35").Value = "'1.00"
$ws.Range("E35").Value = "  +0.08%  "

$ws.Range("D36").Value = "'0.116"
$ws.Range("E36").Value = "  +6.66%  "

$ws.Range("D37").Value = "'157.28"
$ws.Range("E37").Value = "  -3.75%  "

$ws.Range("D38").Value = "'18.75"
$ws.Range("E38").Value = "  -4.14%  "

$ws.Range("D39").Value = "'18.84"
$ws.Range("E39").Value = "  -1.24%  "

$ws.Range("D40").Value = "'1.32"
$ws.Range("E40").Value = "  -4.05%  "

$ws.Range("D41").Value = "'1.00"
$ws.Range("E41").Value = "  +0.06%  "

$ws.Range("D42").Value = "'1.68"
$ws.Range("E42").Value = "  -5.07%  "

$ws.Range("D43").Value = "'4.77"
$ws.Range("E43").Value = "  -4.76%  "

$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").Value = "'2.47"
$ws.Range("E44").Value = "  -3.46%  "

$ws.Range("B45").Value = "PolygonEcosystemToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D45").Value = "'0.321"
$ws.Range("E45").Value = "  -3.47%  "

$ws.Range("D46").Value = "'38.45"
$ws.Range("E46").Value = "  -2.66%  "

$ws.Range("D47").Value = "'145.06"
$ws.Range("E47").Value = "  -7.61%  "

$ws.Range("D48").Value = "'3.56"
$ws.Range("E48").Value = "  -4.23%  "

$ws.Range("D49").Value = "'0.532"
$ws.Range("E49").Value = "  -5.45%  "

$ws.Range("D50").Value = "'1.64"
$ws.Range("E50").Value = "  -5.99%  "

$ws.Range("D51").Value = "'0.597"
$ws.Range("E51").Value = "  -1.77%  "

